$wb = $excel.ActiveWorkbook

function Remove-HyperlinksAt($ws, [string[]]$addrs) {
    $changed = $true
    while ($changed) {
        $changed = $false
        foreach ($hl in $ws.Hyperlinks) {
            $a = $hl.Range.Address()
            if ($addrs -contains $a) {
                $hl.Delete()
                $changed = $true
                break
            }
        }
    }
}

# =====================================================================
# Sheet "Overview": A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
# Insert a new row for cbb64caf-9157-4db8-aa33-dc85b1f6fec4 before the
# existing fb5b9956-379f-4818-ba96-cff1dbe259bf row (pushes it down one).
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows("8:8").Copy()
$wsOverview.Rows("8:8").Insert()

$wsOverview.Range("A8").Value = "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md"
$wsOverview.Range("B8").Value = "Ready for handoff"
$wsOverview.Range("C8").Value = "Ready for handoff"
$wsOverview.Range("D8").Value = "2016-39-20 14:39:15"

Remove-HyperlinksAt $wsOverview @('$A$8')
$wsOverview.Hyperlinks.Add($wsOverview.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/e2e/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md", "", "", "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.md")

# =====================================================================
# Sheet "zh-cn": A=Source File Name, B=File Extension, C=Status,
# D=Latest Handoff File, E=Latest Handoff Datetime, H=Latest Handback
# DateTime, I=Handoff Reason. Same insert pattern as above.
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows("8:8").Copy()
$wsZh.Rows("8:8").Insert()

$wsZh.Range("A8").Value = "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md"
$wsZh.Range("B8").Value = ".md"
$wsZh.Range("C8").Value = "Ready for handoff"
$wsZh.Range("D8").Value = "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.e2e6b4f706b5725c27dcef98bfba3364c8af98e6.zh-cn.xlf"
$wsZh.Range("E8").Value = "2016-03-20 14:39:12"
$wsZh.Range("H8").Value = "0001-01-01 00:00:00"
$wsZh.Range("I8").Value = "Include"

Remove-HyperlinksAt $wsZh @('$A$8','$B$8','$D$8')

$wsZh.Hyperlinks.Add($wsZh.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/e2e/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md", "", "", "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B8"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/e2e/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.e2e6b4f706b5725c27dcef98bfba3364c8af98e6.zh-cn.xlf", "", "", "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.e2e6b4f706b5725c27dcef98bfba3364c8af98e6.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.md")
$wsZh.Hyperlinks.Add($wsZh.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8015b2593e758c2753384d6c443bfda91574c191/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.zh-cn.xlf", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.zh-cn.xlf")

# =====================================================================
# Sheet "de-de": same layout as zh-cn, but for the German target.
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows("8:8").Copy()
$wsDe.Rows("8:8").Insert()

$wsDe.Range("A8").Value = "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md"
$wsDe.Range("B8").Value = ".md"
$wsDe.Range("C8").Value = "Ready for handoff"
$wsDe.Range("D8").Value = "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.e2e6b4f706b5725c27dcef98bfba3364c8af98e6.de-de.xlf"
$wsDe.Range("E8").Value = "2016-03-20 14:39:15"
$wsDe.Range("H8").Value = "0001-01-01 00:00:00"
$wsDe.Range("I8").Value = "Include"

Remove-HyperlinksAt $wsDe @('$A$8','$B$8','$D$8')

$wsDe.Hyperlinks.Add($wsDe.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/e2e/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md", "", "", "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B8"), "https://github.com/OpenLocalizationTest/oltest/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/e2e/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D8"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e2e6b4f706b5725c27dcef98bfba3364c8af98e6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/cbb64caf-9157-4db8-aa33-dc85b1f6fec4.e2e6b4f706b5725c27dcef98bfba3364c8af98e6.de-de.xlf", "", "", "cbb64caf-9157-4db8-aa33-dc85b1f6fec4.e2e6b4f706b5725c27dcef98bfba3364c8af98e6.de-de.xlf")

$wsDe.Hyperlinks.Add($wsDe.Range("A9"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.md")
$wsDe.Hyperlinks.Add($wsDe.Range("B9"), "https://github.com/OpenLocalizationTest/oltest/blob/72f4c1cd46b653c4f64405d84a9aee39bf9ac530/e2e/fb5b9956-379f-4818-ba96-cff1dbe259bf.md", "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D9"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4c4ababb00b10edfcc4bbfa64fb46c998821968c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.de-de.xlf", "", "", "fb5b9956-379f-4818-ba96-cff1dbe259bf.9098578b847812f099eddee2cfa549aae22e7add.de-de.xlf")

Write-Host "Done: inserted cbb64caf-9157-4db8-aa33-dc85b1f6fec4 record into Overview, zh-cn, de-de sheets."
